{"js": "// LOB1235 course sheet update:\n//  - Cr\u00e9ditos-trabalho: 1 -> 0\n//  - Carga hor\u00e1ria: 90 h -> 60 h\n//  - Ativa\u00e7\u00e3o: 01/01/2022 -> 01/01/2025\n//  - Append a sentence about didactic trips to the end of the \"Programa\"\n//    paragraph, in both Portuguese and English.\n\nconst body = context.document.body;\n\n// 1) Cr\u00e9ditos-trabalho: 1 -> 0\nlet results = body.search(\"Cr\u00e9ditos-trabalho: 1\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const item of results.items) {\n  item.insertText(\"Cr\u00e9ditos-trabalho: 0\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Carga hor\u00e1ria: 90 h -> 60 h\nresults = body.search(\"Carga hor\u00e1ria: 90 h\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const item of results.items) {\n  item.insertText(\"Carga hor\u00e1ria: 60 h\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Ativa\u00e7\u00e3o: 01/01/2022 -> 01/01/2025\nresults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2022\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const item of results.items) {\n  item.insertText(\"Ativa\u00e7\u00e3o: 01/01/2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Append sentence to the end of the Portuguese \"Programa\" paragraph\n//    (the one ending in \"Estudos de caso.\").\nresults = body.search(\"Estudos de caso.\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const item of results.items) {\n  item.insertText(\n    \" A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\",\n    Word.InsertLocation.end\n  );\n}\nawait context.sync();\n\n// 5) Append sentence to the end of the English \"Programa\" paragraph\n//    (the one ending in \"Case studies.\").\nresults = body.search(\"Case studies.\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nfor (const item of results.items) {\n  item.insertText(\n    \" The discipline may have didactic trips to complement the content of the discipline.\",\n    Word.InsertLocation.end\n  );\n}\nawait context.sync();\n", "ps1": "# LOB1235 course sheet update:\n#  - Cr\u00e9ditos-trabalho: 1 -> 0\n#  - Carga hor\u00e1ria: 90 h -> 60 h\n#  - Ativa\u00e7\u00e3o: 01/01/2022 -> 01/01/2025\n#  - Append a sentence about didactic trips to the end of the \"Programa\"\n#    paragraph, in both Portuguese and English.\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2\n$wdReplaceAll = 2\n\n# 1) Cr\u00e9ditos-trabalho: 1 -> 0\n$r1 = $d.Range()\n$r1.Find.Execute(\"Cr\u00e9ditos-trabalho: 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"Cr\u00e9ditos-trabalho: 0\", $wdReplaceAll) | Out-Null\n\n# 2) Carga hor\u00e1ria: 90 h -> 60 h\n$r2 = $d.Range()\n$r2.Find.Execute(\"Carga hor\u00e1ria: 90 h\", $false, $false, $false, $false, $false, $true, 1, $false, \"Carga hor\u00e1ria: 60 h\", $wdReplaceAll) | Out-Null\n\n# 3) Ativa\u00e7\u00e3o: 01/01/2022 -> 01/01/2025\n$r3 = $d.Range()\n$r3.Find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2025\", $wdReplaceAll) | Out-Null\n\n# 4) Append sentence to the end of the Portuguese \"Programa\" paragraph\n#    (the one ending in \"Estudos de caso.\").\n$r4 = $d.Range()\n$r4.Find.Execute(\"Estudos de caso.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Estudos de caso. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\", $wdReplaceAll) | Out-Null\n\n# 5) Append sentence to the end of the English \"Programa\" paragraph\n#    (the one ending in \"Case studies.\").\n$r5 = $d.Range()\n$r5.Find.Execute(\"Case studies.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Case studies. The discipline may have didactic trips to complement the content of the discipline.\", $wdReplaceAll) | Out-Null\n"}
